$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginPage")

# Row 3: result text corrected "Invalid" -> "InValid"
$ws.Range("C3").Value = "InValid"

# Row 4: password changed "admin@123" -> "admin#123"; result corrected "Invalid" -> "InValid"
$ws.Range("B4").Value = "admin#123"
$ws.Range("C4").Value = "InValid"

# Move the hyperlink that used to sit on B4 onto B3 instead, keeping the
# same mailto target but showing "admin@123" as its display text.
$ws.Range("B4").Hyperlinks.Delete()
$h = $ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin@123", [Type]::Missing, [Type]::Missing, "admin@123")

# Adding the hyperlink also overwrote B3's cell text with the display text;
# put the real value ("admin123") back while keeping the hyperlink/display.
$ws.Range("B3").Value = "admin123"
$ws.Range("B3").Style = "Hyperlink"

# Update the active cell selection to C5.
$ws.Range("C5").Select()
